$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the first two employees (6 rows: 2 workers x 3 "Periodo Mora"
#    each) from the detail table. The remaining 5 data rows shift up and
#    keep their original content, matching the target layout (rows 16-20).
# ---------------------------------------------------------------------------
$ws.Range("B16:J21").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2) Update the summary figures at the top of the statement.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 498976
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 3

# ---------------------------------------------------------------------------
# 3) Nudge the logo image slightly to the left (same vertical position and
#    size), as was done by hand in the source workbook.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 750450 / 12700.0
$shp.Width = 975600 / 12700.0
$shp.Top = 246450 / 12700.0
$shp.Height = 612000 / 12700.0
